$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency market data refresh (prices + 1h volume deltas).
# A new "OKB" entry is inserted at row 11, which shifts the
# Coin/Link/Price/Volume columns (B:E) down by one row from there on;
# the rank column (A) is left untouched, matching the source diff.
# (The lowest two rows, Cronos/WEMIXTOKEN, fall off the bottom of the
# tracked A1:E51 range as a result of that shift.)

# Pre-format the editable columns as Text so that numeric-looking
# values (e.g. "1.001") are stored as literal strings, matching the
# original inlineStr cell type, instead of being auto-converted to numbers.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '30.500.01'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '2.015.87'
$ws.Range('E3').Value = '  +5.15%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '324.70'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.5115'
$ws.Range('E8').Value = '  +3.10%  '
$ws.Range('D9').Value = '0.08786'
$ws.Range('E9').Value = '  +5.71%  '
$ws.Range('D10').Value = '1.137'
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '43.26'
$ws.Range('E11').Value = '  +2.81%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '24.67'
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '2.012.64'
$ws.Range('E13').Value = '  +4.99%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '6.616'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '7.498'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').Value = '1.008'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '94.54'
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.00001118'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.06526'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '18.93'
$ws.Range('E20').Value = '  +4.26%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.229'
$ws.Range('E22').Value = '  +4.67%  '
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').Value = '30.564.36'
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '11.88'
$ws.Range('E24').Value = '  +5.04%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.229'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.244.86'
$ws.Range('E26').Value = '  +5.04%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '22.43'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '163.28'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.437'
$ws.Range('E29').Value = '  +6.24%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '131.57'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.140'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.1054'
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '6.097'
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.830'
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.358'
$ws.Range('E35').Value = '  +13.64%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.02526'
$ws.Range('E36').Value = '  +3.15%  '
$ws.Range('D37').Value = '5.469'
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06681'
$ws.Range('E38').Value = '  +4.78%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '12.31'
$ws.Range('E39').Value = '  +7.65%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '9.096'
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.2201'
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6687'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.231'
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '13.70'
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.6193'
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '2.199'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '3.667'
$ws.Range('E48').Value = '  +0.71%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.271'
$ws.Range('E49').Value = '  +4.66%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '124.82'
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '81.32'
$ws.Range('E51').Value = '  +3.32%  '

# Restore the default (General) style/number format now that the
# text values are safely stored, so the saved cells carry no extra
# style index, matching the original workbook formatting.
$ws.Range('B2:E51').Style = 'Normal'
